# Resolved bug with Providers: PROVIDER_IDENTIFIER values in column C had
# stray hyphens from a copy/paste (e.g. "GS-07F-BA385" -> should be
# "GS07FBA385"). Strip the hyphens, and also correct the cell formatting
# that tagged along with the bad paste (re-apply the format from the
# neighboring column so the affected cells match the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains("-")) {
        $cell.Value = $val.Replace("-", "")

        # Re-apply the correct (non copy/paste-corrupted) formatting by
        # pulling it from the neighboring FILE_IDENTIFIER cell on the same
        # row, which always carries the correct style.
        $ws.Cells.Item($r, 4).Copy()
        $cell.PasteSpecial(-4122)
    }
}
$excel.CutCopyMode = 0

# Adjust column widths: split the combined G:H width definition so that
# column H and column I get their own (wider) widths.
$ws.Columns.Item(8).ColumnWidth = 28.5
$ws.Columns.Item(9).ColumnWidth = 24.5
